$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "60.442.56"
Set-TextValue "E2" "  +3.37%  "

Set-TextValue "D3" "2.638.94"
Set-TextValue "E3" "  +0.55%  "

Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.01%  "

Set-TextValue "D5" "569.64"
Set-TextValue "E5" "  +6.56%  "

Set-TextValue "D6" "146.52"
Set-TextValue "E6" "  +2.68%  "

Set-TextValue "E7" "  -0.36%  "

Set-TextValue "D8" "0.610"
Set-TextValue "E8" "  +7.65%  "

Set-TextValue "E9" "  -2.13%  "

Set-TextValue "E10" "  +3.96%  "

Set-TextValue "E11" "  +6.30%  "

Set-TextValue "D12" "0.342"
Set-TextValue "E12" "  +2.62%  "

Set-TextValue "D13" "3.106.67"

Set-TextValue "D14" "60.415.64"
Set-TextValue "E14" "  +3.41%  "

Set-TextValue "D15" "21.72"
Set-TextValue "E15" "  +4.72%  "

Set-TextValue "D16" "2.660.46"
Set-TextValue "E16" "  -0.21%  "

Set-TextValue "E17" "  +3.31%  "

Set-TextValue "D18" "4.56"
Set-TextValue "E18" "  +3.80%  "

Set-TextValue "D19" "344.90"
Set-TextValue "E19" "  +3.12%  "

Set-TextValue "D20" "10.42"
Set-TextValue "E20" "  +2.79%  "

Set-TextValue "E21" "  +2.46%  "

Set-TextValue "E22" "  +0.89%  "

Set-TextValue "D23" "0.999"
Set-TextValue "E23" "  +0.07%  "

Set-TextValue "D24" "66.84"
Set-TextValue "E24" "  +1.11%  "

Set-TextValue "D25" "0.443"
Set-TextValue "E25" "  +6.79%  "

Set-TextValue "E26" "  +2.07%  "

Set-TextValue "D27" "0.995"
Set-TextValue "E27" "  -0.44%  "

Set-TextValue "E28" "  +3.34%  "

Set-TextValue "D29" "0.0₃0775"
Set-TextValue "E29" "  +5.52%  "

Set-TextValue "E30" "  -0.14%  "

Set-TextValue "E31" "  +4.63%  "

Set-TextValue "D32" "6.11"
Set-TextValue "E32" "  +4.42%  "

Set-TextValue "D33" "156.31"
Set-TextValue "E33" "  +3.93%  "

Set-TextValue "D34" "19.20"
Set-TextValue "E34" "  +2.53%  "

Set-TextValue "D35" "4.10"
Set-TextValue "E35" "  +5.28%  "

Set-TextValue "B36" "Fetch.AI"
Set-TextValue "C36" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D36" "0.911"
Set-TextValue "E36" "  +12.50%  "

Set-TextValue "B37" "SuiNetwork"
Set-TextValue "C37" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D37" "0.905"
Set-TextValue "E37" "  +6.69%  "

Set-TextValue "E38" "  +6.27%  "

Set-TextValue "D39" "37.59"
Set-TextValue "E39" "  +1.07%  "

Set-TextValue "D40" "1.52"
Set-TextValue "E40" "  +7.56%  "

Set-TextValue "D41" "307.71"
Set-TextValue "E41" "  +9.86%  "

Set-TextValue "D42" "3.67"
Set-TextValue "E42" "  +2.88%  "

Set-TextValue "D43" "0.993"
Set-TextValue "E43" "  -0.59%  "

Set-TextValue "D44" "0.607"
Set-TextValue "E44" "  +1.23%  "

Set-TextValue "D45" "0.0981"
Set-TextValue "E45" "  +4.92%  "

Set-TextValue "D46" "0.0548"
Set-TextValue "E46" "  +3.46%  "

Set-TextValue "D47" "19.45"
Set-TextValue "E47" "  +2.09%  "

Set-TextValue "E48" "  -0.20%  "

Set-TextValue "B49" "VeChain"
Set-TextValue "C49" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D49" "0.0236"
Set-TextValue "E49" "  +5.25%  "

Set-TextValue "B50" "Aave"
Set-TextValue "C50" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D50" "125.88"
Set-TextValue "E50" "  +11.85%  "

Set-TextValue "D51" "1.973.16"
Set-TextValue "E51" "  +1.43%  "
